$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.069.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.453.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.408"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.041.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.450.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.052.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.584.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.183"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.487.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0766"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.52%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.579.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.41%  "
$ws.Range("E51").Value = "  -0.09%  "
